$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 2181.25
$ws.Range("I46").Value = 1430
$ws.Range("J46").Value = 2632
$ws.Range("K46").Value = 4290
$ws.Range("L46").Value = 7896
$ws.Range("M46").Value = -4171
$ws.Range("N46").Value = -8134
$ws.Range("H60").Value = 2181.25
$ws.Range("I60").Value = 1430
$ws.Range("J60").Value = 2632
$ws.Range("K60").Value = 4290
$ws.Range("L60").Value = 7896
$ws.Range("M60").Value = -3806
$ws.Range("N60").Value = -8864
$ws.Range("H80").Value = 2418.516
$ws.Range("I80").Value = 966.0714
$ws.Range("J80").Value = 3614.647
$ws.Range("K80").Value = 2898.2142
$ws.Range("L80").Value = 10843.941
$ws.Range("M80").Value = -1900.2142
$ws.Range("N80").Value = -12839.941
$ws.Range("H83").Value = 2418.516
$ws.Range("I83").Value = 966.0714
$ws.Range("J83").Value = 3614.647
$ws.Range("K83").Value = 8694.642600000001
$ws.Range("L83").Value = 32531.823
$ws.Range("M83").Value = -3702.642600000001
$ws.Range("N83").Value = -42515.823
$ws.Range("H125").Value = 1777.3334
$ws.Range("J125").Value = 2150
$ws.Range("L125").Value = 19350
$ws.Range("N125").Value = -24270
$ws.Range("H132").Value = 1348.5555
$ws.Range("I132").Value = 865.0909
$ws.Range("J132").Value = 6666.6665
$ws.Range("K132").Value = 2595.2727
$ws.Range("L132").Value = 19999.9995
$ws.Range("M132").Value = -65.27269999999999
$ws.Range("N132").Value = -25059.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 320660.88
$ws.Range("I32").Value = 2483.4814
$ws.Range("J32").Value = 3184257.5
$ws.Range("K32").Value = 2483.4814
$ws.Range("L32").Value = 3184257.5
$ws.Range("M32").Value = -2196.4814
$ws.Range("N32").Value = -3184831.5
$ws.Range("H61").Value = 902.1754
$ws.Range("I61").Value = 676.1111
$ws.Range("J61").Value = 1749.9166
$ws.Range("K61").Value = 676.1111
$ws.Range("L61").Value = 1749.9166
$ws.Range("M61").Value = -464.1111
$ws.Range("N61").Value = -2173.9166
$ws.Range("H128").Value = 85999.664
$ws.Range("J128").Value = 85999.664
$ws.Range("L128").Value = 85999.664
$ws.Range("N128").Value = -95959.664
$ws.Range("H132").Value = 957.0909
$ws.Range("I132").Value = 744.95654
$ws.Range("J132").Value = 2041.3334
$ws.Range("K132").Value = 2234.86962
$ws.Range("L132").Value = 6124.0002
$ws.Range("M132").Value = 295.1303800000001
$ws.Range("N132").Value = -11184.0002
$ws.Range("H136").Value = 902.1754
$ws.Range("I136").Value = 676.1111
$ws.Range("J136").Value = 1749.9166
$ws.Range("K136").Value = 2028.3333
$ws.Range("L136").Value = 5249.7498
$ws.Range("M136").Value = 521.6667000000002
$ws.Range("N136").Value = -10349.7498

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3579.8484
$ws.Range("I134").Value = 1156.7317
$ws.Range("J134").Value = 7553.76
$ws.Range("K134").Value = 3470.1951
$ws.Range("L134").Value = 22661.28
$ws.Range("M134").Value = -935.1950999999999
$ws.Range("N134").Value = -27731.28

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7045161.5
$ws.Range("I31").Value = 9435614
$ws.Range("J31").Value = 6607.278
$ws.Range("K31").Value = 9435614
$ws.Range("L31").Value = 6607.278
$ws.Range("M31").Value = -9435319
$ws.Range("N31").Value = -7197.278
$ws.Range("H34").Value = 7045161.5
$ws.Range("I34").Value = 9435614
$ws.Range("J34").Value = 6607.278
$ws.Range("K34").Value = 9435614
$ws.Range("L34").Value = 6607.278
$ws.Range("M34").Value = -9435412
$ws.Range("N34").Value = -7011.278
$ws.Range("H58").Value = 970.9434
$ws.Range("I58").Value = 780.7646999999999
$ws.Range("J58").Value = 1311.2632
$ws.Range("K58").Value = 780.7646999999999
$ws.Range("L58").Value = 1311.2632
$ws.Range("M58").Value = -577.7646999999999
$ws.Range("N58").Value = -1717.2632
$ws.Range("H132").Value = 916.5325
$ws.Range("I132").Value = 852.3538
$ws.Range("J132").Value = 1264.1666
$ws.Range("K132").Value = 2557.0614
$ws.Range("L132").Value = 3792.4998
$ws.Range("M132").Value = -27.06140000000005
$ws.Range("N132").Value = -8852.4998
$ws.Range("H134").Value = 1036.6061
$ws.Range("I134").Value = 978.16364
$ws.Range("J134").Value = 1328.8182
$ws.Range("K134").Value = 2934.49092
$ws.Range("L134").Value = 3986.4546
$ws.Range("M134").Value = -399.4909200000002
$ws.Range("N134").Value = -9056.454600000001
$ws.Range("H136").Value = 970.9434
$ws.Range("I136").Value = 780.7646999999999
$ws.Range("J136").Value = 1311.2632
$ws.Range("K136").Value = 2342.2941
$ws.Range("L136").Value = 3933.7896
$ws.Range("M136").Value = 207.7058999999999
$ws.Range("N136").Value = -9033.7896

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 870.1795
$ws.Range("I5").Value = 567.36365
$ws.Range("J5").Value = 1262.0588
$ws.Range("K5").Value = 1702.09095
$ws.Range("L5").Value = 3786.1764
$ws.Range("M5").Value = -1590.09095
$ws.Range("N5").Value = -4010.1764
$ws.Range("H113").Value = 853.6786
$ws.Range("I113").Value = 594.6875
$ws.Range("J113").Value = 914.6177
$ws.Range("K113").Value = 1784.0625
$ws.Range("L113").Value = 2743.8531
$ws.Range("M113").Value = 385.9375
$ws.Range("N113").Value = -7083.8531
$ws.Range("H135").Value = 870.1795
$ws.Range("I135").Value = 567.36365
$ws.Range("J135").Value = 1262.0588
$ws.Range("K135").Value = 5106.27285
$ws.Range("L135").Value = 11358.5292
$ws.Range("M135").Value = -2571.27285
$ws.Range("N135").Value = -16428.5292

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2041.7931
$ws.Range("I122").Value = 1907
$ws.Range("K122").Value = 5721
$ws.Range("M122").Value = -3271
$ws.Range("H126").Value = 9260814
$ws.Range("I126").Value = 1332.25
$ws.Range("K126").Value = 3996.75
$ws.Range("M126").Value = -1526.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1681.579
$ws.Range("I68").Value = 1652.381
$ws.Range("J68").Value = 1717.6471
$ws.Range("K68").Value = 1652.381
$ws.Range("L68").Value = 1717.6471
$ws.Range("M68").Value = -903.3810000000001
$ws.Range("N68").Value = -3215.6471
$ws.Range("H71").Value = 1681.579
$ws.Range("I71").Value = 1652.381
$ws.Range("J71").Value = 1717.6471
$ws.Range("K71").Value = 8261.905000000001
$ws.Range("L71").Value = 8588.235499999999
$ws.Range("M71").Value = -4517.905000000001
$ws.Range("N71").Value = -16076.2355
$ws.Range("H100").Value = 2647.2903
$ws.Range("I100").Value = 2037.9166
$ws.Range("J100").Value = 3032.158
$ws.Range("K100").Value = 2037.9166
$ws.Range("L100").Value = 3032.158
$ws.Range("M100").Value = -1496.9166
$ws.Range("N100").Value = -4114.157999999999
$ws.Range("H136").Value = 2478.712
$ws.Range("I136").Value = 1941.4222
$ws.Range("J136").Value = 4205.7144
$ws.Range("K136").Value = 5824.2666
$ws.Range("L136").Value = 12617.1432
$ws.Range("M136").Value = -3274.2666
$ws.Range("N136").Value = -17717.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4616
$ws.Range("I62").Value = 4562.75
$ws.Range("J62").Value = 4701.2
$ws.Range("K62").Value = 4562.75
$ws.Range("L62").Value = 4701.2
$ws.Range("M62").Value = -3938.75
$ws.Range("N62").Value = -5949.2
$ws.Range("H65").Value = 4616
$ws.Range("I65").Value = 4562.75
$ws.Range("J65").Value = 4701.2
$ws.Range("K65").Value = 22813.75
$ws.Range("L65").Value = 23506
$ws.Range("M65").Value = -19693.75
$ws.Range("N65").Value = -29746
$ws.Range("H132").Value = 15244716
$ws.Range("I132").Value = 16448124
$ws.Range("J132").Value = 1550
$ws.Range("K132").Value = 49344372
$ws.Range("L132").Value = 4650
$ws.Range("M132").Value = -49341842
$ws.Range("N132").Value = -9710
$ws.Range("H136").Value = 915.5333000000001
$ws.Range("I136").Value = 942.9032
$ws.Range("J136").Value = 854.9286
$ws.Range("K136").Value = 2828.7096
$ws.Range("L136").Value = 2564.7858
$ws.Range("M136").Value = -278.7096000000001
$ws.Range("N136").Value = -7664.7858
